$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = "('Воин Лакватуса', ['{4}{B}{B}', 'Существо — Кошмарный ужас', 'Когда Воин Лакватуса входит в игру, целевой игрок теряет 6 жизней.', 'Когда Воин Лакватуса покидает игру, целевой игрок получает 6 жизней.', '{B}: Восстановить Воина Лакватуса.', '6/3'])"

$ws.Range("A2").Value = $newText
$ws.Range("A3:A22").EntireRow.Delete()
